# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the combined "全部类型" sheet, matching the latest scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (F column holds the count) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 547
$wsExhibit.Range("F4").Value = 195
$wsExhibit.Range("F6").Value = 501
$wsExhibit.Range("F8").Value = 117
$wsExhibit.Range("F10").Value = 6697
$wsExhibit.Range("F12").Value = 367
$wsExhibit.Range("F13").Value = 2968
$wsExhibit.Range("F14").Value = 192
$wsExhibit.Range("F15").Value = 335
$wsExhibit.Range("F17").Value = 537

# --- Sheet "全部类型" (same events, offset by the extra 演出 rows) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 547
$wsAll.Range("F6").Value = 195
$wsAll.Range("F8").Value = 501
$wsAll.Range("F10").Value = 117
$wsAll.Range("F13").Value = 6697
$wsAll.Range("F16").Value = 367
$wsAll.Range("F17").Value = 2968
$wsAll.Range("F18").Value = 192
$wsAll.Range("F19").Value = 335
$wsAll.Range("F21").Value = 537
